$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Isi Burndown Chart" - fill in the Actual Remaining Effort for 5/6 (column O)
# with the same remaining-effort values carried over from 4/6 (column N),
# mirroring the other day-over-day entries already present in the sheet.
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("O13").Value = 2
$ws.Range("O14").Value = 4

# "Edit Back Log" - the reviewer comment on F4 is no longer needed, remove it
$ws.Range("F4").Comment.Delete()

# Leave the selection where data entry ended
[void]$ws.Range("O15").Select()
